# Auto-generated edit script: reshapes Artfynd sheet rows 3-10 per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 112164561
$ws.Range("B3").Value = 93158
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 2667
$ws.Range("F3").Value = 'Platt fjädermossa'
$ws.Range("G3").Value = 'Neckera complanata'
$ws.Range("H3").Value = '(Hedw.) Huebener'
$ws.Range("P3").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q3").Value = 332934.9485370842
$ws.Range("R3").Value = 6626957.391457222
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = 'Värmland'
$ws.Range("U3").Value = 'Eda'
$ws.Range("V3").Value = 'Värmland'
$ws.Range("W3").Value = 'Järnskog'
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = '2022-06-07'
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = '00:00'
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = '2022-06-07'
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = '00:00'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = 'Jeanette Fahlstad'
$ws.Range("AX3").Value = 'Jeanette Fahlstad'

# Row 4
$ws.Range("A4").Value = 112164673
$ws.Range("B4").Value = 93157
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 1078
$ws.Range("F4").Value = 'Rundfjädermossa'
$ws.Range("G4").Value = 'Neckera besseri'
$ws.Range("H4").Value = '(Lobarz.) Jur.'
$ws.Range("P4").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q4").Value = 332854.0101354558
$ws.Range("R4").Value = 6626967.584723449
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = 'Värmland'
$ws.Range("U4").Value = 'Eda'
$ws.Range("V4").Value = 'Värmland'
$ws.Range("W4").Value = 'Järnskog'
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '2022-06-07'
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = '00:00'
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = '2022-06-07'
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = '00:00'
$ws.Range("AC4").Value = 'Under överhängande klippa'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = 'Jeanette Fahlstad'
$ws.Range("AX4").Value = 'Jeanette Fahlstad'

# Row 5
$ws.Range("A5").Value = 112164702
$ws.Range("B5").Value = 89369
$ws.Range("C5").Value = 'Ovaliderad'
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 5447
$ws.Range("F5").Value = 'Vedticka'
$ws.Range("G5").Value = 'Fuscoporia viticola'
$ws.Range("H5").Value = '(Schwein.) Murrill'
$ws.Range("P5").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q5").Value = 332979.8007009666
$ws.Range("R5").Value = 6627033.102525626
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 'Värmland'
$ws.Range("U5").Value = 'Eda'
$ws.Range("V5").Value = 'Värmland'
$ws.Range("W5").Value = 'Järnskog'
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = '2022-06-07'
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = '00:00'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = '2022-06-07'
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = '00:00'
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = 'Jeanette Fahlstad'
$ws.Range("AX5").Value = 'Jeanette Fahlstad'

# Row 6
$ws.Range("A6").Value = 112164607
$ws.Range("B6").Value = 93158
$ws.Range("C6").Value = 'Ovaliderad'
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 2667
$ws.Range("F6").Value = 'Platt fjädermossa'
$ws.Range("G6").Value = 'Neckera complanata'
$ws.Range("H6").Value = '(Hedw.) Huebener'
$ws.Range("P6").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q6").Value = 332973.0533703604
$ws.Range("R6").Value = 6627006.656504014
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 'Värmland'
$ws.Range("U6").Value = 'Eda'
$ws.Range("V6").Value = 'Värmland'
$ws.Range("W6").Value = 'Järnskog'
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = '2022-06-07'
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = '00:00'
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = '2022-06-07'
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = '00:00'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = 'Jeanette Fahlstad'
$ws.Range("AX6").Value = 'Jeanette Fahlstad'

# Row 7
$ws.Range("A7").Value = 112164579
$ws.Range("B7").Value = 93159
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 2666
$ws.Range("F7").Value = 'Grov fjädermossa'
$ws.Range("G7").Value = 'Neckera crispa'
$ws.Range("H7").Value = 'Hedw.'
$ws.Range("P7").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q7").Value = 332922.7263719498
$ws.Range("R7").Value = 6626955.416314425
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Värmland'
$ws.Range("U7").Value = 'Eda'
$ws.Range("V7").Value = 'Värmland'
$ws.Range("W7").Value = 'Järnskog'
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = '2022-06-07'
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = '00:00'
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = '2022-06-07'
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = '00:00'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'Jeanette Fahlstad'
$ws.Range("AX7").Value = 'Jeanette Fahlstad'

# Row 8
$ws.Range("A8").Value = 112164609
$ws.Range("B8").Value = 92683
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 2362
$ws.Range("F8").Value = 'Blek stjärnmossa'
$ws.Range("G8").Value = 'Mnium stellare'
$ws.Range("H8").Value = 'Hedw.'
$ws.Range("P8").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q8").Value = 332973.0533703604
$ws.Range("R8").Value = 6627006.656504014
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Värmland'
$ws.Range("U8").Value = 'Eda'
$ws.Range("V8").Value = 'Värmland'
$ws.Range("W8").Value = 'Järnskog'
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = '2022-06-07'
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = '00:00'
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = '2022-06-07'
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = '00:00'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = 'Jeanette Fahlstad'
$ws.Range("AX8").Value = 'Jeanette Fahlstad'

# Row 9
$ws.Range("A9").Value = 112164661
$ws.Range("B9").Value = 89864
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 5467
$ws.Range("F9").Value = 'Kådvaxskinn'
$ws.Range("G9").Value = 'Phlebia serialis'
$ws.Range("H9").Value = '(Fr.:Fr.) Donk'
$ws.Range("P9").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q9").Value = 332864.8090984516
$ws.Range("R9").Value = 6626971.642313651
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Värmland'
$ws.Range("U9").Value = 'Eda'
$ws.Range("V9").Value = 'Värmland'
$ws.Range("W9").Value = 'Järnskog'
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2022-06-07'
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = '00:00'
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2022-06-07'
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = '00:00'
$ws.Range("AC9").Value = 'På granlåga'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'Jeanette Fahlstad'
$ws.Range("AX9").Value = 'Jeanette Fahlstad'

# Row 10
$ws.Range("A10").Value = 112189008
$ws.Range("B10").Value = 95233
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 2609
$ws.Range("F10").Value = 'Blåsfliksmossa'
$ws.Range("G10").Value = 'Lejeunea cavifolia'
$ws.Range("H10").Value = '(Ehrh.) Lindb.'
$ws.Range("P10").Value = 'Norr Masetjärnet, Vrm'
$ws.Range("Q10").Value = 332854.0101354558
$ws.Range("R10").Value = 6626967.584723449
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Värmland'
$ws.Range("U10").Value = 'Eda'
$ws.Range("V10").Value = 'Värmland'
$ws.Range("W10").Value = 'Järnskog'
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = '2022-06-07'
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = '00:00'
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = '2022-06-07'
$ws.Range("AB10").NumberFormat = "@"
$ws.Range("AB10").Value = '00:00'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'Jeanette Fahlstad'
$ws.Range("AX10").Value = 'Jeanette Fahlstad'

